# InputList/Acoustic/InputList.xlsx
#
# Rename the two worksheets and switch the active tab from
# "Input + Equipment" (now "Channels") to "Snake + PA" (now "Inputs").
#
# Renaming the sheets automatically re-points every formula that refers
# to them by name (e.g. ='Input + Equipment'!C1 -> =Channels!C1), and the
# "Date Updated" text on both sheets recalculates via the existing
# TODAY() formula, so no further action is needed for those.

$wb = $excel.ActiveWorkbook

$channels = $wb.Worksheets.Item("Input + Equipment")
$channels.Name = "Channels"

$inputs = $wb.Worksheets.Item("Snake + PA")
$inputs.Name = "Inputs"

# Make "Inputs" the active/selected sheet (it was "Channels" before).
$inputs.Activate()
